# "imrpoved the statement venn diagram"
#
# 1) Re-cache the (now stale) datetimeFigureOut fields on every layout +
#    the slide master: 11/7/2018 -> 11/23/2018
# 2) Rework the Venn-diagram group on slide 1:
#    - reflow/reposition the boxes (wider boxes, new spacing)
#    - rename Selection's body text ("If, switch" -> "If, else, switch")
#    - rename/expand the While box into an "Iteration" box with
#      While / do / for bullets
#    - drop the standalone "For:" box (folded into Iteration)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder re-cache across slide master + every layout
# ---------------------------------------------------------------------
function Set-DatePlaceholderText {
    param($shapes)
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $tr = $sh.TextFrame.TextRange
            $tr.Text = "ZZZDATEZZZ"
            $tr.Text = "11/23/2018"
        }
    }
}

$design = $p.Designs.Item(1)
$master = $design.SlideMaster
Set-DatePlaceholderText $master.Shapes
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Set-DatePlaceholderText $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Venn diagram rework
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$group = $slide.Shapes.Item(1)

# Ungroup so each rectangle can be repositioned independently.
$ungrouped = $group.Ungroup()

# Shapes, in their original (pre-edit) order:
#   1 Rectangle 3   - big background rectangle
#   2 Rectangle 4   - "Empty:"
#   3 Rectangle 5   - "Declaration:"
#   4 Rectangle 6   - "Expression:"
#   5 TextBox 7     - title "Types of Statements in C++"
#   6 Rectangle 8   - "Selection:"
#   7 Rectangle 10  - "While:"  -> becomes "Iteration:"
#   8 Rectangle 11  - "Block/compound:"
#   9 Rectangle 12  - "For:"    -> removed (merged into Iteration)

$rectBg         = $slide.Shapes.Item(1)
$rectEmpty      = $slide.Shapes.Item(2)
$rectDecl       = $slide.Shapes.Item(3)
$rectExpr       = $slide.Shapes.Item(4)
$titleBox       = $slide.Shapes.Item(5)
$rectSelection  = $slide.Shapes.Item(6)
$rectIteration  = $slide.Shapes.Item(7)
$rectBlock      = $slide.Shapes.Item(8)
$rectFor        = $slide.Shapes.Item(9)

# -- reposition / resize ------------------------------------------------
$rectBg.Left     = 92.07685089111328
$rectBg.Top      = 70.19134521484375
$rectBg.Width    = 834.2308349609375
$rectBg.Height   = 336.0000915527344

$rectEmpty.Left   = 376.5963134765625
$rectEmpty.Top    = 104.99134063720703
$rectEmpty.Width  = 129.60000610351562
$rectEmpty.Height = 266.4000244140625

$rectDecl.Left   = 105.98795318603516
$rectDecl.Top    = 104.99134063720703
$rectDecl.Width  = 129.60000610351562
$rectDecl.Height = 266.4000244140625

$rectExpr.Left   = 241.83677673339844
$rectExpr.Top    = 104.99134063720703
$rectExpr.Width  = 129.60000610351562
$rectExpr.Height = 266.4000244140625

$titleBox.Left   = 382.677978515625
$titleBox.Top    = 376.05181884765625

$rectSelection.Left   = 511.3558349609375
$rectSelection.Top    = 104.99134063720703
$rectSelection.Width  = 129.60000610351562
$rectSelection.Height = 266.4000244140625

$rectIteration.Left   = 647.20458984375
$rectIteration.Top    = 104.99134063720703
$rectIteration.Width  = 129.60000610351562
$rectIteration.Height = 266.4000244140625

$rectBlock.Left   = 783.0534057617188
$rectBlock.Top    = 104.99134063720703
$rectBlock.Width  = 129.60000610351562
$rectBlock.Height = 266.4000244140625

# -- text updates --------------------------------------------------------

# Selection: "If, switch" -> "If, else, switch"
$selTr = $rectSelection.TextFrame.TextRange
$selBody = $selTr.Paragraphs(2, 1)
$selBody.Text = "ZZZSELZZZ"
$selTr.Paragraphs(2, 1).Text = "If, else, switch"

# While -> Iteration, with additional "do," / "for" lines
$iterTr = $rectIteration.TextFrame.TextRange

$iterP1 = $iterTr.Paragraphs(1, 1)
$iterP1.Text = "ZZZITER1ZZZ"
$iterTr.Paragraphs(1, 1).Text = "Iteration:"

$iterP2 = $iterTr.Paragraphs(2, 1)
$iterP2.Text = "ZZZITER2ZZZ"
$iterTr.Paragraphs(2, 1).Text = "While,"

$iterTr.Paragraphs(2, 1).InsertAfter("`rdo,`rfor")

# -- remove the now-redundant standalone "For:" box ----------------------
$rectFor.Delete()

# -- regroup --------------------------------------------------------------
$newRange = $slide.Shapes.Range(@(1, 2, 3, 4, 5, 6, 7, 8))
$newGroup = $newRange.Group()
